# "Generate Report for Handback" — c738f36e was handed back (in sync with
# en-US). Update the localization-status report:
#   - Overview: refresh the status/date columns for both rows
#   - zh-cn / de-de: promote c738f36e to a "handed back" record (new Latest
#     Target File / Latest Handback File / Latest Handback DateTime) and push
#     it above 09d9d419, which keeps its original "Ready for handoff" data

$wb = $excel.ActiveWorkbook

function Clear-AllHyperlinks($ws) {
    while ($ws.Hyperlinks.Count -gt 0) {
        foreach ($hl in $ws.Hyperlinks) {
            $hl.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "09d9d419-10e2-40db-80f1-a03ea1e41a76.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-28-20 00:28:39"

$ov.Range("A3").Value = "c738f36e-d1e0-4800-adfe-ef5cbad142d7.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-28-20 00:28:39"

Clear-AllHyperlinks $ov

$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8186d62ccf7150a9fe64afdde09ed3cc954bd983/e2e/09d9d419-10e2-40db-80f1-a03ea1e41a76.md", "", "", "c738f36e-d1e0-4800-adfe-ef5cbad142d7.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/e50b823b04aad90ae9a81917c98283aa61cd9fc3/e2e/c738f36e-d1e0-4800-adfe-ef5cbad142d7.md", "", "", "09d9d419-10e2-40db-80f1-a03ea1e41a76.md") | Out-Null

# ---------------------------------------------------------------------
# Per-locale sheets ("zh-cn", "de-de")
# ---------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn";
       XlfExt = "zh-cn.xlf";
       HandoffDate2 = "2016-03-20 00:28:52"; HandbackDate2 = "2016-03-20 00:29:13";
       HandoffDate3 = "2016-03-20 00:28:37";
       TargetXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c9cf898d490e96dd1e8bdf8f8bcf7abd453bda23/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c738f36e-d1e0-4800-adfe-ef5cbad142d7.7bb04417058798f7cabaf1f5017023dbdbf327cc.zh-cn.xlf";
       TargetXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/795c348b1d284dc2af595b278676a86b8fa9ea01/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/09d9d419-10e2-40db-80f1-a03ea1e41a76.aaaf0d3495086f877f4c46cae9136dee6d26a8e3.zh-cn.xlf";
     },
    @{ Sheet = "de-de";
       XlfExt = "de-de.xlf";
       HandoffDate2 = "2016-03-20 00:28:55"; HandbackDate2 = "2016-03-20 00:29:18";
       HandoffDate3 = "2016-03-20 00:28:39";
       TargetXlfUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b3cf295bfa546a1c5aa816884e95f0563260fc5f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c738f36e-d1e0-4800-adfe-ef5cbad142d7.7bb04417058798f7cabaf1f5017023dbdbf327cc.de-de.xlf";
       TargetXlfUrl3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2163d9acd0f504aa293c244c9d22458b0992282/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/09d9d419-10e2-40db-80f1-a03ea1e41a76.aaaf0d3495086f877f4c46cae9136dee6d26a8e3.de-de.xlf";
     }
)

$mdUrlC738 = "https://github.com/OpenLocalizationTest/oltest/blob/e50b823b04aad90ae9a81917c98283aa61cd9fc3/e2e/c738f36e-d1e0-4800-adfe-ef5cbad142d7.md"
$mdUrl09d9 = "https://github.com/OpenLocalizationTest/oltest/blob/8186d62ccf7150a9fe64afdde09ed3cc954bd983/e2e/09d9d419-10e2-40db-80f1-a03ea1e41a76.md"

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    $c738Md  = "c738f36e-d1e0-4800-adfe-ef5cbad142d7.md"
    $c738Xlf = "c738f36e-d1e0-4800-adfe-ef5cbad142d7.7bb04417058798f7cabaf1f5017023dbdbf327cc." + $loc.XlfExt
    $d9d9Md  = "09d9d419-10e2-40db-80f1-a03ea1e41a76.md"
    $d9d9Xlf = "09d9d419-10e2-40db-80f1-a03ea1e41a76.aaaf0d3495086f877f4c46cae9136dee6d26a8e3." + $loc.XlfExt

    # Row 2 now holds the c738f36e record: handed back, with new target /
    # handback file + datetime columns populated.
    $ws.Range("A2").Value = $c738Md
    $ws.Range("B2").Value = ".md"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("D2").Value = $c738Xlf
    $ws.Range("E2").Value = $loc.HandoffDate2
    $ws.Range("F2").Value = $c738Md
    $ws.Range("G2").Value = $c738Xlf
    $ws.Range("H2").Value = $loc.HandbackDate2
    $ws.Range("I2").Value = "Include"

    # Row 3 now holds the 09d9d419 record: unchanged "Ready for handoff".
    $ws.Range("A3").Value = $d9d9Md
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $d9d9Xlf
    $ws.Range("E3").Value = $loc.HandoffDate3
    $ws.Range("F3").ClearContents()
    $ws.Range("G3").ClearContents()
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"

    Clear-AllHyperlinks $ws

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrlC738, "", "", $c738Md) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B2"), $mdUrlC738, "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D2"), $loc.TargetXlfUrl2, "", "", $c738Xlf) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrlC738, "", "", $c738Md) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $loc.TargetXlfUrl2, "", "", $c738Xlf) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl09d9, "", "", $d9d9Md) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B3"), $mdUrl09d9, "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D3"), $loc.TargetXlfUrl3, "", "", $d9d9Xlf) | Out-Null

    $ws.Range("F2:G2").Style = "HyperLink"
}

Write-Host "Report regenerated for handback."
